$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(11007, "COXA C/ SOBRECOXA FRANGO CONGELADO CX. 20KG", 46),
    @(11015, "FILE DE PEITO FRANGO CONGELADO CX. 20KG", 32),
    @(11016, "FILE DE PEITO FRANGO CONGELADO BANDEJA", 1),
    @(11023, "DORSO FRANGO CONGELADO CX. 17KG", 12),
    @(11046, "FRANGO CONGELADO CX. 20KG", 7),
    @(11068, "FILE DE PEITO FRANGO CONG. INTERFOLHADO CX. 20KG", 32),
    @(11072, "SOBRECOXA FRANGO RESFRIADA CX. 20KG", 1),
    @(11076, "PEITO FRANGO CONGELADO INDIVIDUAL CX. 20 KG", 3),
    @(11084, "FILE DE SASSAMI FRANGO CONGELADO INTERFOLHADO CX. 20KG", 5),
    @(11085, "CMS CONGELADO CX 20KG", 20),
    @(5023, "BACON DE PALETA", 30),
    @(2003, "picanha", 800),
    @(5015, "test", 56),
    @(11074, "test", 400)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("C15").Select()
